# Apply the "new extraction parameters" update described by the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Rename the sheet (the data came from a fresh extraction run, the
#    sheet is no longer called after the old "Nestle" dataset).
# ---------------------------------------------------------------------
$ws.Name = "Sheet1"

# ---------------------------------------------------------------------
# 2) Update the two mismatched data points (row 9 -> patient id 8) with
#    the values produced by the new extraction parameters.
# ---------------------------------------------------------------------
$ws.Range("B9").Value = 23.31666666666667
$ws.Range("C9").Value = 0.56092948338929882

# ---------------------------------------------------------------------
# 3) Apply the same bordered / bold / centered style used by the header
#    row to the "Patient ID" column (A2:A73), by copying its formats.
# ---------------------------------------------------------------------
$ws.Range("B1").Copy()
$ws.Range("A2:A73").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 4) Column A no longer needs its manual "best fit" width -- put it back
#    to the sheet's standard column width.
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = $ws.StandardWidth

# ---------------------------------------------------------------------
# 5) Refresh the view: scroll further down and select L53 like the
#    author did while reviewing the refreshed numbers.
# ---------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 47
$ws.Range("L53").Select()
